$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2025-05-12 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-05-13 Tuesday", 2)

# Update the division problems in the table (row, column addressed directly
# since several cells share identical text but need different replacements)
$t = $d.Tables(1)

$t.Cell(1,1).Range.Text  = "30÷3=10, 0"
$t.Cell(1,2).Range.Text  = "94÷4=23, 2"
$t.Cell(1,3).Range.Text  = "13÷4=3, 1"
$t.Cell(1,4).Range.Text  = "95÷4=23, 3"
$t.Cell(1,5).Range.Text  = "79÷7=11, 2"

$t.Cell(5,1).Range.Text  = "70÷8=8, 6"
$t.Cell(5,2).Range.Text  = "14÷2=7, 0"
$t.Cell(5,3).Range.Text  = "17÷3=5, 2"
$t.Cell(5,4).Range.Text  = "47÷7=6, 5"
$t.Cell(5,5).Range.Text  = "26÷6=4, 2"

$t.Cell(9,1).Range.Text  = "56÷8=7, 0"
$t.Cell(9,2).Range.Text  = "60÷6=10, 0"
$t.Cell(9,3).Range.Text  = "78÷7=11, 1"
$t.Cell(9,4).Range.Text  = "83÷4=20, 3"
$t.Cell(9,5).Range.Text  = "75÷7=10, 5"

$t.Cell(13,1).Range.Text = "88÷5=17, 3"
$t.Cell(13,2).Range.Text = "63÷6=10, 3"
$t.Cell(13,3).Range.Text = "62÷4=15, 2"
$t.Cell(13,4).Range.Text = "18÷9=2, 0"
$t.Cell(13,5).Range.Text = "27÷5=5, 2"

$t.Cell(17,1).Range.Text = "84÷8=10, 4"
$t.Cell(17,2).Range.Text = "13÷5=2, 3"
$t.Cell(17,3).Range.Text = "25÷2=12, 1"
$t.Cell(17,4).Range.Text = "56÷6=9, 2"
$t.Cell(17,5).Range.Text = "30÷5=6, 0"
